# Auto-generated Excel COM-interop script
# Applies the "scheduled runner" market-data refresh to the Leve Profits workbook.
# For each (sheet, row) this sets H..N (currentAveragePrice.. LeveProfitHQ) to the
# refreshed values. Some rows gain/lose their M (LeveProfitNQ) or N (LeveProfitHQ)
# cell entirely (HQ price now/previously unavailable) -- those use ClearContents()
# so the cell disappears from the sheet rather than being written as 0/blank.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 6
$ws.Range("H6").Value = 3730.6428
$ws.Range("I6").Value = 6378.625
$ws.Range("J6").Value = 200
$ws.Range("K6").Value = 19135.875
$ws.Range("L6").Value = 600
$ws.Range("M6").Value = -19023.875
$ws.Range("N6").Value = -824
# ALC row 62
$ws.Range("H62").Value = 3348.3333
$ws.Range("I62").Value = 2246.6667
$ws.Range("J62").Value = 4450
$ws.Range("K62").Value = 2246.6667
$ws.Range("L62").Value = 4450
$ws.Range("M62").Value = -1622.6667
$ws.Range("N62").Value = -5698
# ALC row 65
$ws.Range("H65").Value = 3348.3333
$ws.Range("I65").Value = 2246.6667
$ws.Range("J65").Value = 4450
$ws.Range("K65").Value = 11233.3335
$ws.Range("L65").Value = 22250
$ws.Range("M65").Value = -8113.333500000001
$ws.Range("N65").Value = -28490
# ALC row 135
$ws.Range("H135").Value = 168.8
$ws.Range("I135").Value = 111
$ws.Range("K135").Value = 999
$ws.Range("M135").Value = 1536
# ALC row 137
$ws.Range("H137").Value = 1346.909
$ws.Range("I137").Value = 1270.8334
$ws.Range("J137").Value = 1438.2
$ws.Range("K137").Value = 3812.5002
$ws.Range("L137").Value = 4314.6
$ws.Range("M137").Value = -1262.5002
$ws.Range("N137").Value = -9414.6
# ALC row 141
$ws.Range("H141").Value = 10933.091
$ws.Range("I141").Value = 12539.444
$ws.Range("K141").Value = 37618.33199999999
$ws.Range("M141").Value = -32438.33199999999

$ws = $wb.Worksheets.Item("ARM")
# ARM row 32
$ws.Range("H32").Value = 4161.946
$ws.Range("I32").Value = 4359.6284
$ws.Range("J32").Value = 702.5
$ws.Range("K32").Value = 4359.6284
$ws.Range("L32").Value = 702.5
$ws.Range("M32").Value = -4072.6284
$ws.Range("N32").Value = -1276.5
# ARM row 45
$ws.Range("H45").Value = 1940.6154
$ws.Range("I45").Value = 1929.8182
$ws.Range("J45").Value = 2000
$ws.Range("K45").Value = 1929.8182
$ws.Range("L45").Value = 2000
$ws.Range("M45").Value = -1552.8182
$ws.Range("N45").Value = -2754
# ARM row 61
$ws.Range("H61").Value = 1924.25
$ws.Range("I61").Value = 1296
$ws.Range("K61").Value = 1296
$ws.Range("M61").Value = -1084
# ARM row 74
$ws.Range("H74").Value = 754.2308
$ws.Range("I74").Value = 754.2308
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 754.2308
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = 119.7692
$ws.Range("N74").ClearContents()
# ARM row 77
$ws.Range("H77").Value = 754.2308
$ws.Range("I77").Value = 754.2308
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 3771.154
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = 596.8459999999995
$ws.Range("N77").ClearContents()
# ARM row 110
$ws.Range("H110").Value = 1083.1892
$ws.Range("I110").Value = 928.4
$ws.Range("J110").Value = 1746.5714
$ws.Range("K110").Value = 928.4
$ws.Range("L110").Value = 1746.5714
$ws.Range("M110").Value = 1116.6
$ws.Range("N110").Value = -5836.5714
# ARM row 122
$ws.Range("H122").Value = 1627.6
$ws.Range("I122").Value = 1374.6666
$ws.Range("J122").Value = 2007
$ws.Range("K122").Value = 4123.9998
$ws.Range("L122").Value = 6021
$ws.Range("M122").Value = -1673.9998
$ws.Range("N122").Value = -10921
# ARM row 136
$ws.Range("H136").Value = 1924.25
$ws.Range("I136").Value = 1296
$ws.Range("K136").Value = 3888
$ws.Range("M136").Value = -1338

$ws = $wb.Worksheets.Item("BSM")
# BSM row 132
$ws.Range("H132").Value = 15500
$ws.Range("J132").Value = 15500
$ws.Range("L132").Value = 15500
$ws.Range("N132").Value = -25620
# BSM row 134
$ws.Range("H134").Value = 16517.285
$ws.Range("I134").Value = 2405.5
$ws.Range("K134").Value = 7216.5
$ws.Range("M134").Value = -4681.5

$ws = $wb.Worksheets.Item("CRP")
# CRP row 7
$ws.Range("H7").Value = 362.1111
$ws.Range("I7").Value = 401.5
$ws.Range("K7").Value = 401.5
$ws.Range("M7").Value = -288.5
# CRP row 22
$ws.Range("H22").Value = 427.77777
$ws.Range("I22").Value = 367.5
$ws.Range("J22").Value = 476
$ws.Range("K22").Value = 367.5
$ws.Range("L22").Value = 476
$ws.Range("M22").Value = -17.5
$ws.Range("N22").Value = -1176
# CRP row 31
$ws.Range("H31").Value = 868.64813
$ws.Range("I31").Value = 727.1395
$ws.Range("J31").Value = 1421.8182
$ws.Range("K31").Value = 727.1395
$ws.Range("L31").Value = 1421.8182
$ws.Range("M31").Value = -432.1395
$ws.Range("N31").Value = -2011.8182
# CRP row 34
$ws.Range("H34").Value = 868.64813
$ws.Range("I34").Value = 727.1395
$ws.Range("J34").Value = 1421.8182
$ws.Range("K34").Value = 727.1395
$ws.Range("L34").Value = 1421.8182
$ws.Range("M34").Value = -525.1395
$ws.Range("N34").Value = -1825.8182
# CRP row 122
$ws.Range("H122").Value = 825.2
$ws.Range("I122").Value = 637.3333
$ws.Range("J122").Value = 1107
$ws.Range("K122").Value = 1911.9999
$ws.Range("L122").Value = 3321
$ws.Range("M122").Value = 538.0001
$ws.Range("N122").Value = -8221

$ws = $wb.Worksheets.Item("CUL")
# CUL row 4
$ws.Range("H4").Value = 409198.62
$ws.Range("I4").Value = 58978.65
$ws.Range("J4").Value = 739961.9399999999
$ws.Range("K4").Value = 176935.95
$ws.Range("L4").Value = 2219885.82
$ws.Range("M4").Value = -176823.95
$ws.Range("N4").Value = -2220109.82
# CUL row 98
$ws.Range("H98").Value = 714
$ws.Range("J98").Value = 1595
$ws.Range("L98").Value = 4785
$ws.Range("N98").Value = -7781
# CUL row 107
$ws.Range("H107").Value = 5521.7144
$ws.Range("I107").Value = 620.5833
$ws.Range("J107").Value = 12056.556
$ws.Range("K107").Value = 1861.7499
$ws.Range("L107").Value = 36169.66800000001
$ws.Range("M107").Value = 58.25009999999997
$ws.Range("N107").Value = -40009.66800000001
# CUL row 113
$ws.Range("H113").Value = 484
$ws.Range("I113").Value = 250.375
$ws.Range("J113").Value = 593.94116
$ws.Range("K113").Value = 751.125
$ws.Range("L113").Value = 1781.82348
$ws.Range("M113").Value = 1418.875
$ws.Range("N113").Value = -6121.82348
# CUL row 131
$ws.Range("H131").Value = 28573000
$ws.Range("I131").Value = 166667440
$ws.Range("J131").Value = 1737.6897
$ws.Range("K131").Value = 500002320
$ws.Range("L131").Value = 5213.0691
$ws.Range("M131").Value = -499997280
$ws.Range("N131").Value = -15293.0691
# CUL row 137
$ws.Range("H137").Value = 9824.608
$ws.Range("I137").Value = 3375
$ws.Range("J137").Value = 11182.421
$ws.Range("K137").Value = 10125
$ws.Range("L137").Value = 33547.263
$ws.Range("M137").Value = -5025
$ws.Range("N137").Value = -43747.263

$ws = $wb.Worksheets.Item("GSM")
# GSM row 126
$ws.Range("H126").Value = 2162.5334
$ws.Range("I126").Value = 1683.8
$ws.Range("K126").Value = 5051.4
$ws.Range("M126").Value = -2581.4

$ws = $wb.Worksheets.Item("LTW")
# LTW row 136
$ws.Range("H136").Value = 2429.75
$ws.Range("I136").Value = 2819.125
$ws.Range("J136").Value = 1651
$ws.Range("K136").Value = 8457.375
$ws.Range("L136").Value = 4953
$ws.Range("M136").Value = -5907.375
$ws.Range("N136").Value = -10053
# LTW row 140
$ws.Range("H140").Value = 45000
$ws.Range("J140").Value = 45000
$ws.Range("L140").Value = 45000
$ws.Range("N140").Value = -55360

$ws = $wb.Worksheets.Item("WVR")
# WVR row 44
$ws.Range("H44").Value = 7694.3335
$ws.Range("I44").Value = 3000
$ws.Range("J44").Value = 10041.5
$ws.Range("K44").Value = 3000
$ws.Range("L44").Value = 10041.5
$ws.Range("M44").Value = -2446
$ws.Range("N44").Value = -11149.5
# WVR row 132
$ws.Range("H132").Value = 2194.8845
$ws.Range("I132").Value = 2098
$ws.Range("J132").Value = 2377.889
$ws.Range("K132").Value = 6294
$ws.Range("L132").Value = 7133.667
$ws.Range("M132").Value = -3764
$ws.Range("N132").Value = -12193.667
# WVR row 136
$ws.Range("H136").Value = 2099.5
$ws.Range("I136").Value = 1866
$ws.Range("J136").Value = 2333
$ws.Range("K136").Value = 5598
$ws.Range("L136").Value = 6999
$ws.Range("M136").Value = -3048
$ws.Range("N136").Value = -12099

Write-Output "Applied scheduled-runner price refresh to 33 leve rows across 8 sheets."
